# Fruta / hortaliza, semanal
# Insert a new weekly record row right before the current row 443 in the
# "Femacal de La Calera - Pepino ensalada" data table, shifting all the
# following rows (old 443..471) down by one (they become 444..472).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 443; everything below (incl. the old
# row 443) shifts down one row, and the used range grows to A1:R472.
$ws.Rows.Item(443).EntireRow.Insert()

# Populate the newly inserted row 443 with the new weekly observation.
$ws.Cells.Item(443, 1).Value = 3
$ws.Cells.Item(443, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(443, 3).Value = "Coquimbo"
$ws.Cells.Item(443, 4).Value = 44931
$ws.Cells.Item(443, 5).Value = 5
$ws.Cells.Item(443, 6).Value = 100112043
$ws.Cells.Item(443, 7).Value = "Pepino ensalada"
$ws.Cells.Item(443, 8).Value = "Sin especificar"
$ws.Cells.Item(443, 9).Value = "Primera"
$ws.Cells.Item(443, 10).Value = 95
$ws.Cells.Item(443, 11).Value = 18000
$ws.Cells.Item(443, 12).Value = 19000
$ws.Cells.Item(443, 13).Value = 18526
$ws.Cells.Item(443, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(443, 15).Value = "Limache"
$ws.Cells.Item(443, 16).Value = 265
$ws.Cells.Item(443, 17).Value = 70
$ws.Cells.Item(443, 18).Value = "Hortaliza"
